$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.260.78"
$ws.Range("E2").Value = "  -1.26%  "
$ws.Range("D3").Value = "2.770.85"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "353.34"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").Value = "107.94"
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("E7").Value = "  -3.33%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "0.583"
$ws.Range("E9").Value = "  -1.23%  "
$ws.Range("D10").Value = "39.53"
$ws.Range("E10").Value = "  -1.30%  "
$ws.Range("E11").Value = "  +3.29%  "
$ws.Range("D12").Value = "20.03"
$ws.Range("E12").Value = "  +3.44%  "
$ws.Range("E13").Value = "  -2.23%  "
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("D15").Value = "3.204.79"
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").Value = "2.769.40"
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").Value = "51.196.08"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("D19").Value = "7.62"
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("D20").Value = "3.08"
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("D21").Value = "13.08"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").Value = "0.0₃0959"
$ws.Range("E22").Value = "  -1.60%  "
$ws.Range("D23").Value = "69.66"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("D24").Value = "265.36"
$ws.Range("E24").Value = "  -3.20%  "
$ws.Range("D25").Value = "2.70"
$ws.Range("E25").Value = "  -1.16%  "
$ws.Range("E26").Value = "  +0.16%  "
$ws.Range("D27").Value = "25.93"
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("D28").Value = "0.162"
$ws.Range("E28").Value = "  +13.18%  "
$ws.Range("D29").Value = "10.17"
$ws.Range("E29").Value = "  +0.49%  "
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("D31").Value = "36.09"
$ws.Range("E31").Value = "  +7.06%  "
$ws.Range("D32").Value = "6.12"
$ws.Range("E32").Value = "  +7.37%  "
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("D34").Value = "0.0440"
$ws.Range("E34").Value = "  -5.39%  "
$ws.Range("D35").Value = "5.46"
$ws.Range("E35").Value = "  +3.86%  "
$ws.Range("E36").Value = "  -1.98%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").Value = "18.08"
$ws.Range("E38").Value = "  +0.33%  "
$ws.Range("D39").Value = "3.14"
$ws.Range("E39").Value = "  -2.56%  "
$ws.Range("E40").Value = "  -1.69%  "
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("E42").Value = "  -0.56%  "
$ws.Range("D43").Value = "120.41"
$ws.Range("E43").Value = "  -0.75%  "
$ws.Range("D44").Value = "22.00"
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("D46").Value = "2.095.50"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "2.32"
$ws.Range("E47").Value = "  +2.44%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "3.24"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("D49").Value = "0.902"
$ws.Range("E49").Value = "  -1.74%  "
$ws.Range("D50").Value = "5.39"
$ws.Range("E50").Value = "  -5.17%  "
$ws.Range("E51").Value = "  +7.25%  "
